$wb = $excel.ActiveWorkbook

# --- shredding sheet -------------------------------------------------
$wsShredding = $wb.Worksheets.Item("shredding")
$wsShredding.Activate() | Out-Null
$wsShredding.Range("G3").Value = 0.207
$wsShredding.Range("G4").Select() | Out-Null

# --- extrusion sheet ---------------------------------------------------
$wsExtrusion = $wb.Worksheets.Item("extrusion")
$wsExtrusion.Activate() | Out-Null
$wsExtrusion.Range("G3").Value = 0.207
$wsExtrusion.Range("G4").Select() | Out-Null

# --- granulate sheet ----------------------------------------------------
$wsGranulate = $wb.Worksheets.Item("granulate")
$wsGranulate.Activate() | Out-Null
$wsGranulate.Range("G3").Value = 0.207
$wsGranulate.Range("G4").Select() | Out-Null

# --- conditioning sheet: ends up the active / selected tab --------------
$wsConditioning = $wb.Worksheets.Item("conditioning")
$wsConditioning.Activate() | Out-Null
$wsConditioning.Range("Q3").Select() | Out-Null
